# Code Merge Changes - 9/25/2017
# Update the "TestResultExcelFilePath" value (column H, row 2) on the
# payroll-processing / report sheets to point at the new shared test-result
# output location.

$wb = $excel.ActiveWorkbook

$newPath = "F:\\Automation_TestResults\\Payroll_Tax_StatutoryScenarios\\Automation Test Result for Statutory Scenarios201718.xlsx"

$targetSheets = @(
    "ProcessPayrollForJulyMonthSAPP",
    "ProcessPayrollForAugMonthSAPP",
    "ProcessPayrollForSepMonthSAPP",
    "AverageWeeklyEarningsTestReport",
    "ProcessPayrollForJan16MonthSAPP"
)

foreach ($name in $targetSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("H2").Value = $newPath
}
